# Swap the order of "Recorded By" names in column G:
# "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# Only cells that contain exactly the old value are touched; cells that
# already say "System" or "dnasr281@gmail.com" alone are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldVal = "dnasr281@gmail.com, System"
$newVal = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G is the 7th column
    if ($cell.Value2 -eq $oldVal) {
        $cell.Value = $newVal
    }
}
